# daily auto push: 2025-10-05 18:34 UTC
# Append the next day's row (2025/10/06, 月, hour 1, rank 6) to the bottom
# of the tracking sheet (row 66), matching the format of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

# Column A holds a date-looking string ("2025/10/06") that must stay plain
# text (like every other row in the sheet) instead of being auto-converted
# to a real Excel date serial number. Temporarily force the cell to Text
# format before assigning it, then clear the formatting again so the cell
# is left with the workbook's normal/default style (same as the rest of
# the sheet) while keeping the text value that was already committed.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/06"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 6
